# Transition doc edit: "Added my email to transition doc"
#
# The Jake Pennington (Developer) bullet is the only stakeholder entry
# that lacks a contact e-mail/hyperlink (every other team member has
# " <hyperlink>" appended after the "(Role):" label). This script adds
# Jake's e-mail as a live mailto hyperlink, matching the pattern used
# for the other names, and relocates the document's "_GoBack" bookmark
# (which Word always re-drops at the most-recently-edited spot) to the
# blank paragraph immediately below, which is where it ends up once the
# edit above becomes the last edit made to the document.

$d = $word.ActiveDocument

# --- locate the "Jake Pennington (Developer):" paragraph -----------------
$targetIndex = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $paraText = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13)
    if ($paraText -eq "Jake Pennington (Developer):") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -ne $null) {

    $emailText    = "jake.pennington@umontana.edu"
    $emailAddress = "mailto:" + $emailText

    # Append a space, then the e-mail text, to the end of that paragraph
    # (same layout as "Mark Matas (Team Lead/Developer): <email>", etc.).
    $para = $d.Paragraphs.Item($targetIndex)
    $para.Range.InsertAfter(" " + $emailText)

    # Re-fetch the paragraph range now that it has grown, and figure out
    # where the plain e-mail text we just inserted sits (Range.End points
    # one past the paragraph mark, hence the "-1").
    $para      = $d.Paragraphs.Item($targetIndex)
    $paraEnd   = $para.Range.End - 1
    $emailStart = $paraEnd - $emailText.Length
    $emailRange = $d.Range($emailStart, $paraEnd)

    # Turn that plain text into a live mailto hyperlink, same as the
    # other stakeholders' e-mail addresses.
    $d.Hyperlinks.Add($emailRange, $emailAddress, $null, $null, $emailText) | Out-Null

    # --- move the "_GoBack" bookmark to the blank paragraph right below ---
    if ($d.Bookmarks.Exists("_GoBack")) {
        $d.Bookmarks.Item("_GoBack").Delete()
    }
    $blankPara = $d.Paragraphs.Item($targetIndex + 1)
    $d.Bookmarks.Add("_GoBack", $blankPara.Range) | Out-Null
}
